# Add the new "Create Truck preview module" task as row 14 of the plan sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A14").Value = "#100012"
$ws.Range("B14").Value = "Create Truck preview module"

# Move the selection to where the author left it after the edit.
$ws.Range("I13").Select()
